# Generate Report for handoff
#
# - Updates the "Handoff transform failed" status text to "Ready for
#   handoff" everywhere it appears (Overview!B2/C2, zh-cn!B2, de-de!B2).
# - Adds a "Latest Handoff File" hyperlink (column C, row 2) on the
#   zh-cn and de-de sheets, pointing at the freshly generated xlf file
#   for each locale.
# - Refreshes the handoff datetime / handback datetime / handoff reason
#   columns for row 2 on both locale sheets to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# 1. Status text change (shared by Overview/zh-cn/de-de "Status" cells)
# ---------------------------------------------------------------------
$wsOverview.Range("B2").Value2 = $newStatus
$wsOverview.Range("C2").Value2 = $newStatus
$wsZhCn.Range("B2").Value2 = $newStatus
$wsDeDe.Range("B2").Value2 = $newStatus

# ---------------------------------------------------------------------
# Helper values
# ---------------------------------------------------------------------
$repoBlobBase = "https://github.com/OpenLocalizationTest/oltest/blob/03e0516f524fa473df9eed9f3c7517dbdfecb66f/e2e/"

$zhCnXlfName = "e763d729-3c0a-492f-a7c8-9d63441ec1c3.2936bdbc8579b0f45d1488d9df20663323d26946.zh-cn.xlf"
$deDeXlfName = "e763d729-3c0a-492f-a7c8-9d63441ec1c3.2936bdbc8579b0f45d1488d9df20663323d26946.de-de.xlf"

$zhCnXlfUrl = $repoBlobBase + $zhCnXlfName
$deDeXlfUrl = $repoBlobBase + $deDeXlfName

$includeText = "Include"
$ignoredText = "Ignored"
$defaultDate = "0001-01-01 00:00:00"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/03e0516f524fa473df9eed9f3c7517dbdfecb66f/e2e/e763d729-3c0a-492f-a7c8-9d63441ec1c3.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/03e0516f524fa473df9eed9f3c7517dbdfecb66f/.localization-config"

# ---------------------------------------------------------------------
# 2. zh-cn sheet - row 2 (handoff just happened)
# ---------------------------------------------------------------------
$wsZhCn.Range("D2").Value2 = "2016-01-26 12:15:11"
$wsZhCn.Range("G2").Value2 = $defaultDate
$wsZhCn.Range("H2").Value2 = $includeText

# row 3 values are unchanged text, but re-assert them so the shared
# strings line up with the (reused) strings above
$wsZhCn.Range("D3").Value2 = $defaultDate
$wsZhCn.Range("G3").Value2 = $defaultDate
$wsZhCn.Range("H3").Value2 = $ignoredText

# Rebuild the hyperlinks collection in top-to-bottom, left-to-right
# order (A2, C2, A3) so the new "Latest Handoff File" link lands
# between the existing two, matching the document's natural reading
# order.
$zhCnA2Display = $wsZhCn.Range("A2").Value2
$zhCnA3Display = $wsZhCn.Range("A3").Value2

$wsZhCn.Hyperlinks.Delete() | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, $zhCnA2Display) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhCnXlfUrl, [Type]::Missing, [Type]::Missing, $zhCnXlfName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, $zhCnA3Display) | Out-Null

# Recreating the hyperlinks re-applies Excel's default (theme-colored)
# hyperlink look to every touched cell; restore the workbook's existing
# custom hyperlink font (underline, Calibri 11, #6495ED) on all of them
# so the visual style matches the rest of the workbook.
foreach ($ref in @("A2", "C2", "A3")) {
    $r = $wsZhCn.Range($ref)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 11
    $r.Font.Underline = 2
    $r.Font.Color = 0xED9564
}

# ---------------------------------------------------------------------
# 3. de-de sheet - row 2 (handoff just happened)
# ---------------------------------------------------------------------
$wsDeDe.Range("D2").Value2 = "2016-01-26 12:15:22"
$wsDeDe.Range("G2").Value2 = $defaultDate
$wsDeDe.Range("H2").Value2 = $includeText

$wsDeDe.Range("D3").Value2 = $defaultDate
$wsDeDe.Range("G3").Value2 = $defaultDate
$wsDeDe.Range("H3").Value2 = $ignoredText

$deDeA2Display = $wsDeDe.Range("A2").Value2
$deDeA3Display = $wsDeDe.Range("A3").Value2

$wsDeDe.Hyperlinks.Delete() | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, $deDeA2Display) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deDeXlfUrl, [Type]::Missing, [Type]::Missing, $deDeXlfName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, $deDeA3Display) | Out-Null

foreach ($ref in @("A2", "C2", "A3")) {
    $r = $wsDeDe.Range($ref)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 11
    $r.Font.Underline = 2
    $r.Font.Color = 0xED9564
}
